$d = $word.ActiveDocument

# --- Step 1 ----------------------------------------------------------------
# Grab a copy of the "Meta description" paragraph (the 2nd paragraph) before
# it gets deleted. We reuse its exact run layout (a leading empty run
# followed by a bold run) for the new heading paragraph we add near the end
# of the document.
$metaP = $d.Paragraphs(2)
$metaP.Range.Copy()

# --- Step 2 ----------------------------------------------------------------
# Insert a brand-new paragraph right before the very last paragraph (the one
# that used to hold the "Create a captivating feature image ..." image
# prompt) and paste the copied run layout into it.
$n = $d.Paragraphs.Count
$prevP = $d.Paragraphs($n - 1)
$prevP.Range.InsertParagraphAfter()
$newP = $d.Paragraphs($n)
$newP.Style = "Normal"
$insertRange = $newP.Range
$insertRange.Collapse(1)
$insertRange.Paste()

# Trim the pasted "Meta description: Discover Fruit Rainbow, ..." text down
# to just the bold heading "Play Fruit Rainbow Slot for Free - Review and
# Ratings 2021", keeping the leading empty run and bold run/formatting.
$targetP = $d.Paragraphs($n)
$labelLen = "Meta description".Length
$pStart = $targetP.Range.Start
$pEnd = $targetP.Range.End
$restRange = $d.Range($pStart + $labelLen, $pEnd - 1)
$restRange.Text = ""
$boldRange = $d.Range($pStart, $pStart + $labelLen)
$boldRange.Text = "Play Fruit Rainbow Slot for Free - Review and Ratings 2021"

# --- Step 3 ----------------------------------------------------------------
# Replace the italic "Create a captivating feature image ..." image-prompt
# text (still the last paragraph) with the meta-description sentence.
$lastP = $d.Paragraphs($d.Paragraphs.Count)
$lastP.Range.Find.Execute(
    "Create a captivating feature image for Fruit Rainbow slot game, featuring a happy Maya warrior with glasses in a cartoon style. The Maya warrior should be positioned in front of a vibrant rainbow and surrounded by various fruit symbols, including strawberries, oranges, cherries, bananas, and lime slices. The background should be bursting with vibrant colors, showcasing a dynamic and lively ambiance that captures the essence of the game. Additionally, the image should vividly depict the stacked wilds and scatters to highlight the excitement of the game. Overall, the feature image should convey the fun nature of this fruit-themed slot machine while emphasizing the opportunity for massive wins.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Discover Fruit Rainbow, a high variance slot game with 40 paylines and stacked wilds. Play for free and read our review and ratings for 2021.",
    2)

# --- Step 4 ----------------------------------------------------------------
# Finally, delete the original "Meta description" paragraph near the top of
# the document (it has now been relocated to the bottom).
$metaP2 = $d.Paragraphs(2)
$metaP2.Range.Delete()
